$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 4-9 (6 rows): shifts everything below up by 6.
$ws.Range("A4:A9").EntireRow.Delete()

# Replace the formula in Q3 with its computed literal value (validation result).
$ws.Range("Q3").Value = 5

# Update the selection to match the post-edit active cell.
$ws.Range("Q3").Select()
